# Kawasaki vehicle config workbook - add the Meta/Datorama daily email subject
# row to the "5525" sheet (mirrors the existing gmail_subject row pattern).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("5525")

# New label/value pair appended as row 9
$ws.Range("A9").Value = "meta_gmail_subject"
$ws.Range("B9").Value = 'Datorama | Report "Kawasaki Daily Reporting Meta" (1159468)'

# Match the formatting that was pasted in alongside the new row: a plain
# black label next to a dark-grey value in a different font face.
$ws.Range("A9").Font.Color = 0

$ws.Range("B9").Font.Name = "Aptos Display"
$ws.Range("B9").Font.Color = 2039583

# Excel leaves the cursor one row below the newly entered data.
$ws.Range("A10").Select()
